$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This is the "Generate Report for Handback" edit: the localization-status
# report is refreshed after a handback. The "Ready for handoff" status
# becomes "Handed back: in sync with en-US" everywhere it appears, and the
# zh-cn / de-de detail sheets get their "Latest Target File", "Latest
# Handback File" and "Latest Handback DateTime" columns (I/J/K) populated
# for both data rows, including a new hyperlink on the "Latest Target File"
# cell pointing at the same source-file URL used by column A.
# ---------------------------------------------------------------------------

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetFile = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5273e89801e7eb3bdb830e1cd62e15c664c05b3d/e2e/b488e0a4-e7d8-4343-a46b-53f4ea708df8.md"
$sourceUrlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5273e89801e7eb3bdb830e1cd62e15c664c05b3d/e2e/ffffb1983d0c-e3f9-4ca6-80eb-d06f1cff7aff.md"
$displayB   = "ffffb1983d0c-e3f9-4ca6-80eb-d06f1cff7aff.md"

# ---------------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) just mirror the status.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn / de-de detail sheets: same shape, per-language handback timestamp
# and xlf file name.
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; XlfFile = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.fa112066d26f78bbcbadb52fba6ff71b07da0b4a.zh-cn.xlf"; HandbackTime = "2016-08-30 11:13:12" },
    @{ Name = "de-de"; XlfFile = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.fa112066d26f78bbcbadb52fba6ff71b07da0b4a.de-de.xlf"; HandbackTime = "2016-08-30 11:13:19" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Handback File (J) / Latest Handback DateTime (K) text for both rows.
    $ws.Range("J2").Value = $lang.XlfFile
    $ws.Range("K2").Value = $lang.HandbackTime
    $ws.Range("J3").Value = $lang.XlfFile
    $ws.Range("K3").Value = $lang.HandbackTime

    # Rebuild the hyperlinks collection so the new "Latest Target File" (I)
    # links land in A2, I2, A3, I3 order, matching the source-file link.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFile)
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFile)
    $ws.Hyperlinks.Add($ws.Range("A3"), $sourceUrlB, [Type]::Missing, [Type]::Missing, $displayB)
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFile)

    # Column widths: Status (C) and the newly-populated Target/Handback File
    # columns (I/J) widen to fit the longer text.
    $ws.Columns.Item(3).ColumnWidth = 29.17
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
